{"js": "const pairs = [\n  [\"2024-12-05 Thursday\", \"2024-12-06 Friday\"],\n  [\"918\u00f74=\", \"752\u00f74=\"],\n  [\"374\u00f73=\", \"147\u00f78=\"],\n  [\"221\u00f79=\", \"628\u00f75=\"],\n  [\"793\u00f76=\", \"776\u00f78=\"],\n  [\"367\u00f72=\", \"820\u00f79=\"],\n  [\"229\u00f73=\", \"554\u00f78=\"],\n  [\"600\u00f76=\", \"401\u00f74=\"],\n  [\"640\u00f76=\", \"619\u00f79=\"],\n  [\"560\u00f76=\", \"818\u00f73=\"],\n  [\"756\u00f78=\", \"258\u00f75=\"],\n  [\"359\u00f79=\", \"687\u00f72=\"],\n  [\"830\u00f73=\", \"694\u00f77=\"],\n  [\"438\u00f73=\", \"166\u00f74=\"],\n  [\"670\u00f73=\", \"605\u00f73=\"],\n  [\"904\u00f77=\", \"626\u00f75=\"],\n  [\"679\u00f75=\", \"974\u00f74=\"],\n  [\"352\u00f72=\", \"221\u00f73=\"],\n  [\"559\u00f75=\", \"182\u00f75=\"],\n  [\"635\u00f73=\", \"974\u00f77=\"],\n  [\"938\u00f78=\", \"702\u00f72=\"],\n  [\"289\u00f77=\", \"826\u00f77=\"],\n  [\"502\u00f72=\", \"532\u00f73=\"],\n  [\"822\u00f79=\", \"872\u00f75=\"],\n  [\"584\u00f73=\", \"587\u00f76=\"],\n  [\"988\u00f74=\", \"289\u00f79=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the twenty-five \"NNN\u00f7N=\" division problems to\n# the new set of values (output generated at c986bee).\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{ Old = \"2024-12-05 Thursday\"; New = \"2024-12-06 Friday\" },\n    @{ Old = \"918\u00f74=\"; New = \"752\u00f74=\" },\n    @{ Old = \"374\u00f73=\"; New = \"147\u00f78=\" },\n    @{ Old = \"221\u00f79=\"; New = \"628\u00f75=\" },\n    @{ Old = \"793\u00f76=\"; New = \"776\u00f78=\" },\n    @{ Old = \"367\u00f72=\"; New = \"820\u00f79=\" },\n    @{ Old = \"229\u00f73=\"; New = \"554\u00f78=\" },\n    @{ Old = \"600\u00f76=\"; New = \"401\u00f74=\" },\n    @{ Old = \"640\u00f76=\"; New = \"619\u00f79=\" },\n    @{ Old = \"560\u00f76=\"; New = \"818\u00f73=\" },\n    @{ Old = \"756\u00f78=\"; New = \"258\u00f75=\" },\n    @{ Old = \"359\u00f79=\"; New = \"687\u00f72=\" },\n    @{ Old = \"830\u00f73=\"; New = \"694\u00f77=\" },\n    @{ Old = \"438\u00f73=\"; New = \"166\u00f74=\" },\n    @{ Old = \"670\u00f73=\"; New = \"605\u00f73=\" },\n    @{ Old = \"904\u00f77=\"; New = \"626\u00f75=\" },\n    @{ Old = \"679\u00f75=\"; New = \"974\u00f74=\" },\n    @{ Old = \"352\u00f72=\"; New = \"221\u00f73=\" },\n    @{ Old = \"559\u00f75=\"; New = \"182\u00f75=\" },\n    @{ Old = \"635\u00f73=\"; New = \"974\u00f77=\" },\n    @{ Old = \"938\u00f78=\"; New = \"702\u00f72=\" },\n    @{ Old = \"289\u00f77=\"; New = \"826\u00f77=\" },\n    @{ Old = \"502\u00f72=\"; New = \"532\u00f73=\" },\n    @{ Old = \"822\u00f79=\"; New = \"872\u00f75=\" },\n    @{ Old = \"584\u00f73=\"; New = \"587\u00f76=\" },\n    @{ Old = \"988\u00f74=\"; New = \"289\u00f79=\" }\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Wrap = 1          # wdFindContinue\n    $find.Execute($pair.Old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null\n}\n"}
